$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: push the existing "Table1" (No / Nama Barang / Jumlah /
#    Keterangan, rows 9:10) down two rows to 11:12, freeing rows 8 and 9 for
#    two new summary rows ("Jumlah Unit" / "Total Harga").
# ---------------------------------------------------------------------------
$ws.Rows("9:10").Insert()

# ---------------------------------------------------------------------------
# 2. Make room for two new table columns ("Harga (Rp.)" and
#    "Sub Total (Rp.)") by inserting two blank columns at E:F - this pushes
#    the existing "Keterangan" column from E to G.
# ---------------------------------------------------------------------------
$ws.Columns("E:F").Insert()

# ---------------------------------------------------------------------------
# 3. Re-sync the table (ListObject) to its new location/extent right away,
#    before any new text is written into the old header area, so the table
#    picks up the existing column headers correctly.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B11:G12"))

# ---------------------------------------------------------------------------
# 4. Fill the two new summary rows (copy the look of the row above first).
# ---------------------------------------------------------------------------
$ws.Range("B7").Copy()
$ws.Range("B8:B9").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C8:C9").PasteSpecial(-4122)

$ws.Range("B8").Value = "Jumlah Unit"
$ws.Range("B9").Value = "Total Harga"
$ws.Range("C8").Value = "[onshow.jumlahunit]"
$ws.Range("C9").Value = "[onshow.totalharga]"

# ---------------------------------------------------------------------------
# 5. Name the table's new/shifted columns.
# ---------------------------------------------------------------------------
$hdr = $lo.HeaderRowRange
$hdr.Cells(1,6).Value = "Keterangan"
$hdr.Cells(1,5).Value = "Sub Total (Rp.)"
$hdr.Cells(1,4).Value = "Harga (Rp.)"

# ---------------------------------------------------------------------------
# 6. Fill the new data-row cells under the two new columns.
# ---------------------------------------------------------------------------
$ws.Range("F12").Value = "[a.subtotal]"
$ws.Range("E12").Value = "[a.harga]"

# ---------------------------------------------------------------------------
# 7. Column widths for the two newly inserted columns.
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 21
$ws.Columns("F").ColumnWidth = 26.666666666666668

# ---------------------------------------------------------------------------
# 8. Selection, matching what Excel left selected after the edit.
# ---------------------------------------------------------------------------
[void]$ws.Range("E12").Select()
